# Auto update on 2025-12-24 15:08:26
# Apply updated values to jama_exports/kp_data.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 91
$ws.Range("O3").Value = 91
$ws.Range("R3").Value = 0.4

# Row 4
$ws.Range("P4").Value = 115
$ws.Range("Q4").Value = 2.95

# Row 6
$ws.Range("J6").Value = 414
$ws.Range("N6").Value = 414
$ws.Range("Q6").Value = 4.99
$ws.Range("R6").Value = 0.14

# Row 7
$ws.Range("P7").Value = 96
$ws.Range("Q7").Value = 4.22

# Row 9
$ws.Range("F9").Value = 195
$ws.Range("N9").Value = 195
$ws.Range("P9").Value = 79
$ws.Range("Q9").Value = 2.47

# Row 10
$ws.Range("F10").Value = 245
$ws.Range("G10").Value = 87
$ws.Range("N10").Value = 245
$ws.Range("O10").Value = 87
$ws.Range("P10").Value = 77

# Row 11
$ws.Range("F11").Value = 50
$ws.Range("N11").Value = 50
$ws.Range("P11").Value = 25
$ws.Range("Q11").Value = 2

# Row 12
$ws.Range("F12").Value = 100
$ws.Range("N12").Value = 337
$ws.Range("P12").Value = 106
$ws.Range("Q12").Value = 3.18

# Row 14
$ws.Range("F14").Value = 320
$ws.Range("N14").Value = 320
$ws.Range("P14").Value = 96
$ws.Range("Q14").Value = 3.33

# Row 15
$ws.Range("F15").Value = 42
$ws.Range("N15").Value = 85
$ws.Range("P15").Value = 28
$ws.Range("R15").Value = 0.09

# Row 16
$ws.Range("F16").Value = 71
$ws.Range("N16").Value = 74
$ws.Range("P16").Value = 44
$ws.Range("Q16").Value = 1.68
$ws.Range("R16").Value = 0.09

# Row 17
$ws.Range("F17").Value = 28
$ws.Range("N17").Value = 28
$ws.Range("P17").Value = 8
$ws.Range("Q17").Value = 3.5
